# "maj emilien probleme fin d'annee"
# Adds 3 new exercise rows to the "exos" inventory sheet:
#   1. A new "fichiers" entry (FIC-013) appended to the existing fichiers block
#      (right after FIC-012, which is the original row 96).
#   2. A new "stationnaire" entry (STATIO-006) appended to the existing
#      stationnaire block (right after STATIO-005, the original row 115).
#   3. A brand new "problemes" entry (puissance velo), appended at the very
#      end of the data, just before the trailing "fin" marker row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert the three new (still empty) rows, top to bottom, so the
#     row numbers below already account for the shift caused by each
#     previous insertion. ---

# New row for "fichiers" / FIC-013, right after the current last "fichiers"
# row (old row 96).
$ws.Rows("97").Insert()

# New row for "stationnaire" / STATIO-006, right after the current last
# "stationnaire" row (old row 115, now at row 116 because of the insert
# above).
$ws.Rows("117").Insert()

# New row for the brand new "problemes" domain, right before the trailing
# "fin" row (old row 128, now at row 130 because of the two inserts above).
$ws.Rows("130").Insert()

# --- Step 2: fill in column A (domain) for the two rows that reuse an
#     already-existing shared string, so these writes don't create new
#     shared-string entries. ---
$ws.Range("A97").Value = "fichiers"
$ws.Range("A117").Value = "stationnaire"

# --- Step 3: fill in the brand-new text values, in the same order they were
#     first typed by the author (this matters only for shared-string
#     table ordering, not for the visible result): STATIO-006 block first,
#     then FIC-013 block, then the new "problemes" row. ---

$ws.Range("B117").Value = "STATIO-006"
$ws.Range("C117").Value = "Méthodes numériques"

$ws.Range("B97").Value = "FIC-013"
$ws.Range("C97").Value = "Lecture d'un fichier donnant les décimales de pi"

$ws.Range("A130").Value = "problemes"
$ws.Range("B130").Value = "puissance vélo"
$ws.Range("C130").Value = "Analyse de la puissance d'un cycliste à partir d'un relevé GPS"

# --- Step 4: restore the selection to match the author's final cursor
#     position (bottom of the newly extended list). ---
$ws.Range("C131").Select()
